$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'72.207.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.26%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'4.040.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.04%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'539.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.10%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'152.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.99%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'4.032.78"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.04%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.698"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.11%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E10").Value = "'  -1.14%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.173"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.02%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'53.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +9.31%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000331"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.32%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -0.18%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.682.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.07%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'4.038.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.01%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'14.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.04%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'20.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.61%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -0.74%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -0.98%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'72.143.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.35%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'448.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +2.79%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'97.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.37%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -1.89%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -0.49%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'14.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.74%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'4.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +15.15%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.24%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'10.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.38%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'5.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.07%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'37.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.38%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'8.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +17.13%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.134"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.80%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'13.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.85%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'49.25"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +14.76%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'680.30"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.33%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'66.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.85%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.453"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +4.44%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.0₃0897"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +5.71%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.148"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -6.08%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'WEMIXToken"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C41").Style = "Normal"
$ws.Range("E41").Value = "'  -4.75%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'ThetaToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'3.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.17%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'THORChain"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'11.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +16.24%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +0.04%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'FirstDigitalUSD"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.09%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'VeChain"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'0.0492"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.37%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -0.65%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.35%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'3.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E50").Value = "'  -3.73%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'3.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.48%  "
$ws.Range("E51").Style = "Normal"
